# TC08_INS_CancerType-GastricCancer.xlsx – INS regression-suite update
#
# The "ProgramsTab" SQL query stored in cell B2 is rewritten: the plain
# `prg.website AS "Website"` column is replaced with a CASE expression
# that falls back between `prg.program_acronym` / `prg.program_link`,
# and the indentation of the existing "Data Location Details" CASE is
# tightened by a few spaces. The cell keeps the same wrap-text / size-12
# look, the sheet's scroll position / selection is nudged down to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the Program-query text held in B2 --------------------------
$newQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$cell = $ws.Range("B2")
$cell.Value = $newQuery

# Re-touch the font (same face/size it already has) so the cell picks up a
# freshly-written style record, matching the re-saved workbook's style churn
# while leaving the visible formatting (wrap text, size 12) untouched.
$cell.Font.Name = $cell.Font.Name

# --- 2. Move the on-screen selection down to B8 ----------------------------
$ws.Range("B8").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
